$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165)
    3  = @(3.272327238179451, 1.626987699542094, 3993.344853322108, 0.5333859586016987, 3998.777554218431)
    4  = @(0.01253208636536152, 0.04103571897497393, 3.223369029078222, 13.86384647080068, 17.14078330521924)
    5  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    6  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    7  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    8  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    9  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    10 = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.104883657715537)
    11 = @(0.2881169905109251, 9.983522426115931, 3.223369029078222, 2797.565817734744, 2811.060826180449)
    12 = @(0.1169995834814548, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 19.67207418041587)
    13 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    14 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    15 = @(0.6545652718822623, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.716211508195562)
    16 = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.104883657715537)
    17 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
